$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column before column I (date), shifting date/legislator_name/legislator_id
# one column to the right. This makes room for the new "category" column and inherits
# the formatting (style) of the surrounding cells automatically.
$ws.Columns.Item(9).Insert()

# Header row: new "category" column (I), plus two appended columns (M=source_file, N=index)
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# Copy formatting from the existing header cell (H1) onto the two newly appended header cells.
$ws.Range($ws.Cells.Item(1, 8), $ws.Cells.Item(1, 8)).Copy()
$ws.Range($ws.Cells.Item(1, 13), $ws.Cells.Item(1, 14)).PasteSpecial(-4122)

# Data rows
$rows = @(2, 3, 4)
foreach ($r in $rows) {
    $idxVal = $ws.Cells.Item($r, 1).Value()

    $ws.Cells.Item($r, 9).Value = "normal"
    $ws.Cells.Item($r, 13).Value = "tmp14431"
    $ws.Cells.Item($r, 14).Value = $idxVal
}

# Copy formatting from the existing data cell (H column) onto the newly appended data columns.
$ws.Range($ws.Cells.Item(2, 8), $ws.Cells.Item(4, 8)).Copy()
$ws.Range($ws.Cells.Item(2, 13), $ws.Cells.Item(4, 14)).PasteSpecial(-4122)

$excel.CutCopyMode = 0
